$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.305.77'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.931.15'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7546'
$ws.Range('E5').Value = '  +5.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.61'
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3175'
$ws.Range('E8').Value = '  -2.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.52'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07003'
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7798'
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07983'
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('D13').Value = '1.931.42'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.360'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.24'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.42'
$ws.Range('E16').Value = '  -2.52%  '
$ws.Range('D17').Value = '30.296.57'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '252.39'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007909'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.729'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = '2.188.28'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.668'
$ws.Range('E24').Value = '  -3.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.496'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.89'
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1335'
$ws.Range('E27').Value = '  +3.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.209'
$ws.Range('E29').Value = '  -5.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.363'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.111'
$ws.Range('E33').Value = '  -2.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05156'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.285'
$ws.Range('E35').Value = '  +1.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7456'
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.771'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01946'
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.793'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '77.45'
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.400'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4459'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.964'
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8303'
$ws.Range('E45').Value = '  -1.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.75'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.718'
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.453'
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '981.59'
$ws.Range('E49').Value = '  +11.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.32'
$ws.Range('E50').Value = '  +1.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06004'
$ws.Range('E51').Value = '  -0.97%  '
